# Apply the "LG and Rrr results summarized; gait metric added to SS-FP file" edit.
#
# This adds a new column U to Sheet1 that records the "gait metric" (unit of
# measurement) used by each cohort/study, plus a header label and an updated
# selection/zoom on the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-study gait metric units (column U, rows 4-12), entered in the same
# order the shared strings first appear in the saved workbook (cm/s, m/s, s,
# then the "gait metric" header last).
$ws.Range("U4").Value  = "cm/s"   # EAS
$ws.Range("U5").Value  = "m/s"    # ELSA
$ws.Range("U6").Value  = "s"      # HRS
$ws.Range("U7").Value  = "s"      # ILSE
$ws.Range("U8").Value  = "s"      # LASA
$ws.Range("U9").Value  = "m/s"    # MAP
$ws.Range("U10").Value = "s"      # NuAge
$ws.Range("U11").Value = "s"      # OCTO
$ws.Range("U12").Value = "s"      # SATSA

# Header for the new column.
$ws.Range("U2").Value = "gait metric"

# Update the view: zoom to 110% and move the active selection to U3.
$win = $excel.ActiveWindow
$win.Zoom = 110
$ws.Range("U3").Select()
